$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.233.70'
$ws.Range("E2").Value = '  +2.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.428.36'
$ws.Range("E3").Value = '  +1.89%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.23'
$ws.Range("E5").Value = '  +1.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.48'
$ws.Range("E6").Value = '  +3.66%  '

$ws.Range("E7").Value = '  +0.95%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.502'
$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.39'
$ws.Range("E10").Value = '  +3.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0802'
$ws.Range("E11").Value = '  +1.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.124'
$ws.Range("E12").Value = '  +2.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.83'
$ws.Range("E13").Value = '  +1.94%  '

$ws.Range("E14").Value = '  +1.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.806.51'
$ws.Range("E15").Value = '  +2.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.431.64'
$ws.Range("E16").Value = '  +2.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.836'
$ws.Range("E17").Value = '  +3.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.195.86'
$ws.Range("E18").Value = '  +2.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.34'
$ws.Range("E19").Value = '  +1.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.44'
$ws.Range("E20").Value = '  +1.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0908'
$ws.Range("E21").Value = '  +2.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.62'
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.81'
$ws.Range("E23").Value = '  +2.37%  '

$ws.Range("E24").Value = '  +1.72%  '

$ws.Range("E25").Value = '  +1.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.27'
$ws.Range("E27").Value = '  +1.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.36'
$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.60'
$ws.Range("E29").Value = '  +4.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.33'
$ws.Range("E30").Value = '  +5.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.117'
$ws.Range("E31").Value = '  +13.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.71'
$ws.Range("E32").Value = '  +8.99%  '

$ws.Range("E33").Value = '  +1.57%  '

$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0763'
$ws.Range("E35").Value = '  +3.25%  '

$ws.Range("E36").Value = '  +3.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.50'
$ws.Range("E37").Value = '  +4.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '129.51'
$ws.Range("E38").Value = '  +24.32%  '

$ws.Range("E39").Value = '  +4.23%  '

$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("E41").Value = '  +0.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.04'
$ws.Range("E42").Value = '  -6.92%  '

$ws.Range("E43").Value = '  +2.71%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.957.79'
$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("E45").Value = '  +1.72%  '

$ws.Range("E46").Value = '  +4.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.43'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.66'
$ws.Range("E48").Value = '  +9.97%  '

$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.49'
$ws.Range("E49").Value = '  +1.36%  '

$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.52'
$ws.Range("E50").Value = '  +2.34%  '

$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.16'
$ws.Range("E51").Value = '  +1.10%  '
